# Apply a different built-in table style to the three tables in the
# deck (slides 14, 15 and 16), matching the commit's table-style change
# from {E1277EE5-020C-4C6A-ACCB-454B1B5BCE5A} to
# {24139B49-052D-4460-BADF-7EA9031C21A0}.

$p = $ppt.ActivePresentation

$newStyleId = "{24139B49-052D-4460-BADF-7EA9031C21A0}"

$slideIndexesWithTables = @(14, 15, 16)

foreach ($idx in $slideIndexesWithTables) {
    $slide = $p.Slides.Item($idx)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
